$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) cells - force text storage to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.255.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.987.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.985.24'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.499'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.08'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.479.91'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.982.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.207.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '431.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0932'
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.58'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '49.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0659'
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.109'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '384.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.631.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.73'
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '30.96'
$ws.Range("D51").Style = "Normal"

# Update other text cells (Coin name, Link, Volume%)
$ws.Range("E2").Value = '  -4.27%  '
$ws.Range("E3").Value = '  -5.73%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  -2.05%  '
$ws.Range("E6").Value = '  -6.73%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -5.62%  '
$ws.Range("E9").Value = '  -2.87%  '
$ws.Range("E10").Value = '  -5.73%  '
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("E13").Value = '  -5.37%  '
$ws.Range("E14").Value = '  -4.98%  '
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("E16").Value = '  -5.78%  '
$ws.Range("E17").Value = '  -5.93%  '
$ws.Range("E18").Value = '  -4.37%  '
$ws.Range("E19").Value = '  -4.60%  '
$ws.Range("E20").Value = '  -5.88%  '
$ws.Range("E21").Value = '  -6.32%  '
$ws.Range("E22").Value = '  -4.67%  '
$ws.Range("E23").Value = '  -7.16%  '
$ws.Range("E24").Value = '  -3.92%  '
$ws.Range("E25").Value = '  -3.86%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("E28").Value = '  -4.08%  '
$ws.Range("E29").Value = '  -4.48%  '
$ws.Range("E30").Value = '  -6.46%  '
$ws.Range("E31").Value = '  -8.44%  '
$ws.Range("E32").Value = '  -6.49%  '
$ws.Range("E33").Value = '  -7.85%  '
$ws.Range("E34").Value = '  -7.68%  '
$ws.Range("E35").Value = '  -7.28%  '
$ws.Range("E36").Value = '  -3.07%  '
$ws.Range("E37").Value = '  -3.04%  '
$ws.Range("E38").Value = '  -6.13%  '
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("E40").Value = '  -6.67%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("E41").Value = '  -1.42%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("E42").Value = '  -4.14%  '
$ws.Range("E43").Value = '  -6.67%  '
$ws.Range("E44").Value = '  -6.45%  '
$ws.Range("E46").Value = '  -5.92%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("E47").Value = '  -5.07%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E48").Value = '  -4.54%  '
$ws.Range("E49").Value = '  -3.77%  '
$ws.Range("E50").Value = '  -6.03%  '
$ws.Range("E51").Value = '  -10.23%  '
